$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 161.66667
$ws.Range("I6").Value = 161.66667
$ws.Range("K6").Value = 485.00001
$ws.Range("M6").Value = -373.00001

$ws.Range("H121").Value = 1978.3334
$ws.Range("J121").Value = 1978.3334
$ws.Range("L121").Value = 5935.0002
$ws.Range("N121").Value = -9429.0002

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = ""

$ws.Range("H138").Value = 2274.3333
$ws.Range("J138").Value = 2697.6924
$ws.Range("L138").Value = 8093.0772
$ws.Range("N138").Value = -18373.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5674.125
$ws.Range("I32").Value = 5752.298
$ws.Range("K32").Value = 5752.298
$ws.Range("M32").Value = -5465.298

$ws.Range("H61").Value = 2577.1
$ws.Range("I61").Value = 2847.125
$ws.Range("J61").Value = 1497
$ws.Range("K61").Value = 2847.125
$ws.Range("L61").Value = 1497
$ws.Range("M61").Value = -2635.125
$ws.Range("N61").Value = -1921

$ws.Range("H74").Value = 3076.8096
$ws.Range("J74").Value = 2499.5
$ws.Range("L74").Value = 2499.5
$ws.Range("N74").Value = -4247.5

$ws.Range("H77").Value = 3076.8096
$ws.Range("J77").Value = 2499.5
$ws.Range("L77").Value = 12497.5
$ws.Range("N77").Value = -21233.5

$ws.Range("H110").Value = 7919.7
$ws.Range("I110").Value = 8899.75
$ws.Range("J110").Value = 3999.5
$ws.Range("K110").Value = 8899.75
$ws.Range("L110").Value = 3999.5
$ws.Range("M110").Value = -6854.75
$ws.Range("N110").Value = -8089.5

$ws.Range("H132").Value = 6553.5835
$ws.Range("I132").Value = 6507.6763
$ws.Range("J132").Value = 7334
$ws.Range("K132").Value = 19523.0289
$ws.Range("L132").Value = 22002
$ws.Range("M132").Value = -16993.0289
$ws.Range("N132").Value = -27062

$ws.Range("H136").Value = 2577.1
$ws.Range("I136").Value = 2847.125
$ws.Range("J136").Value = 1497
$ws.Range("K136").Value = 8541.375
$ws.Range("L136").Value = 4491
$ws.Range("M136").Value = -5991.375
$ws.Range("N136").Value = -9591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 414.85715
$ws.Range("I2").Value = 414.85715
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 414.85715
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -301.85715
$ws.Range("N2").Value = ""

$ws.Range("H5").Value = 750.625
$ws.Range("I5").Value = 350.83334
$ws.Range("J5").Value = 1950
$ws.Range("K5").Value = 350.83334
$ws.Range("L5").Value = 1950
$ws.Range("M5").Value = -238.83334
$ws.Range("N5").Value = -2174

$ws.Range("H31").Value = 1704.2927
$ws.Range("I31").Value = 1605.9333
$ws.Range("J31").Value = 1972.5454
$ws.Range("K31").Value = 1605.9333
$ws.Range("L31").Value = 1972.5454
$ws.Range("M31").Value = -1310.9333
$ws.Range("N31").Value = -2562.5454

$ws.Range("H34").Value = 1704.2927
$ws.Range("I34").Value = 1605.9333
$ws.Range("J34").Value = 1972.5454
$ws.Range("K34").Value = 1605.9333
$ws.Range("L34").Value = 1972.5454
$ws.Range("M34").Value = -1403.9333
$ws.Range("N34").Value = -2376.5454

$ws.Range("H59").Value = 40583.082
$ws.Range("I59").Value = 45400
$ws.Range("J59").Value = 37142.43
$ws.Range("K59").Value = 45400
$ws.Range("L59").Value = 37142.43
$ws.Range("M59").Value = -44255
$ws.Range("N59").Value = -39432.43

$ws.Range("H60").Value = 10553.909
$ws.Range("I60").Value = 10553.909
$ws.Range("K60").Value = 10553.909
$ws.Range("M60").Value = -10042.909

$ws.Range("H99").Value = 1706.4
$ws.Range("I99").Value = 1706.4
$ws.Range("K99").Value = 1706.4
$ws.Range("M99").Value = -208.4000000000001

$ws.Range("H126").Value = 1706.4
$ws.Range("I126").Value = 1706.4
$ws.Range("K126").Value = 5119.200000000001
$ws.Range("M126").Value = -2649.200000000001

$ws.Range("H132").Value = 2166.3914
$ws.Range("I132").Value = 2137
$ws.Range("J132").Value = 2475
$ws.Range("K132").Value = 6411
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -3881
$ws.Range("N132").Value = -12485

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 768618.5600000001
$ws.Range("I4").Value = 19095.975
$ws.Range("K4").Value = 57287.925
$ws.Range("M4").Value = -57175.925

$ws.Range("H132").Value = 2997
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

$ws.Range("H137").Value = 2097.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3392.2
$ws.Range("I22").Value = 3392.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3392.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2863.2
$ws.Range("N22").Value = ""

$ws.Range("H80").Value = 1538.7142
$ws.Range("I80").Value = 1545.3334
$ws.Range("J80").Value = 1499
$ws.Range("K80").Value = 1545.3334
$ws.Range("L80").Value = 1499
$ws.Range("M80").Value = -547.3334
$ws.Range("N80").Value = -3495

$ws.Range("H83").Value = 1538.7142
$ws.Range("I83").Value = 1545.3334
$ws.Range("J83").Value = 1499
$ws.Range("K83").Value = 7726.666999999999
$ws.Range("L83").Value = 7495
$ws.Range("M83").Value = -2734.666999999999
$ws.Range("N83").Value = -17479

$ws.Range("H107").Value = 1916.4286
$ws.Range("J107").Value = 4198.75
$ws.Range("L107").Value = 4198.75
$ws.Range("N107").Value = -8038.75

$ws.Range("H122").Value = 1760.3636
$ws.Range("I122").Value = 995.625
$ws.Range("K122").Value = 2986.875
$ws.Range("M122").Value = -536.875

$ws.Range("H132").Value = 2878.5264
$ws.Range("I132").Value = 2587.1875
$ws.Range("J132").Value = 4432.3335
$ws.Range("K132").Value = 7761.5625
$ws.Range("L132").Value = 13297.0005
$ws.Range("M132").Value = -5231.5625
$ws.Range("N132").Value = -18357.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6614.0435
$ws.Range("I40").Value = 4636.2
$ws.Range("K40").Value = 4636.2
$ws.Range("M40").Value = -4500.2

$ws.Range("H122").Value = 2820.5
$ws.Range("I122").Value = 2820.5
$ws.Range("K122").Value = 8461.5
$ws.Range("M122").Value = -6011.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 3000
$ws.Range("J51").Value = 3000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -4020

$ws.Range("H107").Value = 4895.069
$ws.Range("I107").Value = 4245.2104
$ws.Range("J107").Value = 6129.8
$ws.Range("K107").Value = 12735.6312
$ws.Range("L107").Value = 18389.4
$ws.Range("M107").Value = -10815.6312
$ws.Range("N107").Value = -22229.4

$ws.Range("H126").Value = 1756.9474
$ws.Range("I126").Value = 1666.3334
$ws.Range("J126").Value = 2096.75
$ws.Range("K126").Value = 4999.0002
$ws.Range("L126").Value = 6290.25
$ws.Range("M126").Value = -2529.0002
$ws.Range("N126").Value = -11230.25

